# Generate Report for Handoff
#
# Re-runs the handoff-report generator: the six rows whose "Status" is
# "Ready for handoff" and whose most recent handoff timestamp matches the
# previous report run now get a fresh generation timestamp ("Latest HO
# Xliff Generate Date" on Overview, "Latest Handoff Datetime" on the
# zh-cn / de-de sheets), and on the zh-cn sheet their Priority is marked
# "ht" (handoff type).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-17 04:20:13"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-17 04:20:02"
}

# --- de-de sheet: "Latest Handoff Datetime" (H) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 8).Value = "2016-08-17 04:20:13"
}
